$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - first occurrence of the event rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2821
$ws1.Range("F4").Value = 130

# Sheet "全部类型" (all types) - same events duplicated/aggregated here
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 2821
$ws4.Range("F8").Value = 130
